$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row text updates ---
$ws.Range("B1").Value = "Student Name"
$ws.Range("E1").Value = "Student Tag"

# --- Row 2: Sridhar, Kavitha (Student) ---
$ws.Range("B2").Value = "Sridhar, Kavitha (Student)"
$ws.Range("B2").Borders.LineStyle = -4142
$ws.Range("B2").Font.Color = 2368548
$ws.Range("B2").Font.Name = "Segoe UI"

$ws.Range("C2").Value = "k.sridhar@my.ccsu.edu"
$ws.Range("C2").HorizontalAlignment = -4108
$ws.Range("C2").Borders.LineStyle = 1
$ws.Range("C2").Borders.Weight = -4138
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:k.sridhar@my.ccsu.edu")

# --- Row 3: Singh, Surya P. (Student) ---
$ws.Range("B3").Value = "Singh, Surya P. (Student)"
$ws.Range("B3").Borders.LineStyle = -4142
$ws.Range("B3").Font.Color = 2368548
$ws.Range("B3").Font.Name = "Segoe UI"

$ws.Range("C3").Value = "suryasingh@my.ccsu.edu"
$ws.Range("C3").HorizontalAlignment = -4108
$ws.Range("C3").Borders.LineStyle = 1
$ws.Range("C3").Borders.Weight = -4138

# --- Updated start-time text shared by both rows (appended last to match string table order) ---
$ws.Range("G2").Value = "10/29/2024  22:00:00 PM"
$ws.Range("G3").Value = "10/29/2024  22:00:00 PM"

# --- Remove the last (now-blank) row 17 so used range shrinks to G16 ---
$ws.Rows("17:17").Delete()

# --- Selection moves from G8 to F8 ---
$ws.Range("F8").Select()
